# Fruta / hortaliza, semanal
# Update the "Fecha" (D), "Calidad" (I), "Volumen" (J), "Precio máximo" (L),
# "Precio promedio ponderado" (M) and "Precio $/Kg" (P) columns so the data
# rows reflect the re-sorted weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44642

# Row 3
$ws.Range("D3").Value = 44656
$ws.Range("J3").Value = 100

# Row 4
$ws.Range("D4").Value = 44628
$ws.Range("I4").Value = "Primera"
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861

# Row 5
$ws.Range("D5").Value = 44651

# Row 6
$ws.Range("D6").Value = 44384
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 60
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 833

# Row 7
$ws.Range("D7").Value = 44659
$ws.Range("J7").Value = 80

# Row 8
$ws.Range("D8").Value = 44637
$ws.Range("J8").Value = 100

# Row 9
$ws.Range("D9").Value = 44658
$ws.Range("J9").Value = 80

# Row 10
$ws.Range("D10").Value = 44645

# Row 11
$ws.Range("D11").Value = 44635
$ws.Range("J11").Value = 100

# Row 13
$ws.Range("D13").Value = 44630
$ws.Range("J13").Value = 60

# Row 14
$ws.Range("D14").Value = 44649
